$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "Kaium" -> "Тимур", scores 42/43 -> text "100"/"99"
$ws.Range("A4").Value = "Тимур"

$ws.Range("B4").Value = "'100"
$ws.Range("B4").ClearFormats()

$ws.Range("C4").Value = "'99"
$ws.Range("C4").ClearFormats()
